$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @{
    "G2" = 72.266001
    "H2" = 216.798003
    "I2" = 0.2949652269937106
    "J2" = 0.2949652269937106
    "M2" = 10.67287833333334
    "N2" = 32.018635
    "O2" = 0.2083050184412124
    "P2" = 0.2083050184412124
    "Q2" = 771.2862363095452
    "R2" = 6941.576126785905
    "S2" = 0.06144273704844128
    "T2" = 0.06144273704844127

    "G3" = 72.266001
    "H3" = 216.798003
    "I3" = 0.2949652269937106
    "J3" = 0.2949652269937106
    "O3" = 0.5074067008595954
    "P3" = 0.5074067008595954
    "Q3" = 1878.763207496552
    "R3" = 16908.86886746897
    "S3" = 0.1496673326971804
    "T3" = 0.1496673326971804

    "G4" = 72.266001
    "H4" = 216.798003
    "I4" = 0.2949652269937106
    "J4" = 0.2949652269937106
    "M4" = 14.566016
    "N4" = 43.698048
    "O4" = 0.2842882806991923
    "P4" = 0.2842882806991923
    "Q4" = 1052.627726822016
    "R4" = 9473.649541398143
    "S4" = 0.08385515724808895
    "T4" = 0.08385515724808895

    "H5" = 410.023338
    "I5" = 0.5578585839920717
    "J5" = 0.5578585839920718
    "M5" = 10.67287833333334
    "N5" = 32.018635
    "O5" = 0.2083050184412124
    "P5" = 0.2083050184412124
    "Q5" = 1458.709733433737
    "R5" = 13128.38760090363
    "S5" = 0.1162047426260571
    "T5" = 0.1162047426260571

    "H6" = 410.023338
    "I6" = 0.5578585839920717
    "J6" = 0.5578585839920718
    "O6" = 0.5074067008595954
    "P6" = 0.5074067008595954
    "Q6" = 3553.246575104858
    "S6" = 0.2830611836496226
    "T6" = 0.2830611836496226

    "H7" = 410.023338
    "I7" = 0.5578585839920717
    "J7" = 0.5578585839920718
    "M7" = 14.566016
    "N7" = 43.698048
    "O7" = 0.2842882806991923
    "P7" = 0.2842882806991923
    "Q7" = 1990.802167227136
    "R7" = 17917.21950504422
    "S7" = 0.158592657716392
    "T7" = 0.158592657716392

    "G8" = 36.057927
    "H8" = 108.173781
    "I8" = 0.1471761890142177
    "J8" = 0.1471761890142177
    "M8" = 10.67287833333334
    "N8" = 32.018635
    "O8" = 0.2083050184412124
    "P8" = 0.2083050184412124
    "Q8" = 384.8418678232151
    "R8" = 3463.576810408935
    "S8" = 0.03065753876671398
    "T8" = 0.03065753876671397

    "G9" = 36.057927
    "H9" = 108.173781
    "I9" = 0.1471761890142177
    "J9" = 0.1471761890142177
    "O9" = 0.5074067008595954
    "P9" = 0.5074067008595954
    "Q9" = 937.429851503704
    "R9" = 8436.868663533336
    "S9" = 0.07467818451279243
    "T9" = 0.07467818451279243

    "G10" = 36.057927
    "H10" = 108.173781
    "I10" = 0.1471761890142177
    "J10" = 0.1471761890142177
    "M10" = 14.566016
    "N10" = 43.698048
    "O10" = 0.2842882806991923
    "P10" = 0.2842882806991923
    "Q10" = 525.220341608832
    "R10" = 4726.983074479487
    "S10" = 0.0418404657347113
    "T10" = 0.0418404657347113
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
